$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Column A width change (15.7109375 -> 27.5703125 raw OOXML width)
$ws.Columns.Item(1).ColumnWidth = 26.67

# Row 7: A7 becomes a CONCATENATE formula (was static text "TariffDecision3"),
# inheriting the left/wrap/bordered style already used by C7/D7 (s=13)
$ws.Range("C7").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$ws.Range("A7").Formula = '=CONCATENATE("Tariff ",C7,D7,"-",E7)'
$ws.Rows.Item(7).RowHeight = 90

# Row 8: A8 becomes a CONCATENATE formula (was static text "TariffDecision4")
$ws.Range("C8").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null
$ws.Range("A8").Formula = '=CONCATENATE("Tariff ",C8,D8,"-",E8)'
$ws.Rows.Item(8).RowHeight = 45

$excel.CutCopyMode = $false

# Selection moves from A6:F6 to B8
$ws.Range("B8").Select() | Out-Null
